$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values: B2 and D2 updated, C2 and E2 cleared (removed)
$ws.Range("B2").Value = 28.55068550296312
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 17.824412475279871
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 23.738576534686445
$ws.Range("C3").Value = 23.27312537304482
$ws.Range("D3").Value = 15.569199809149438
$ws.Range("E3").Value = 29.129066284357918

# Update the selected range to reflect the new selection B1:E3
$ws.Range("B1:E3").Select()
